# Updates cryptos list cell values (Price / Volume(1h) columns) per source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.893.51"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").Value = "2.237.17"
$ws.Range("E3").Value = "  -1.80%  "
$ws.Range("E4").Value = "  +0.14%  "
$c = $ws.Range("D5")
$c.Value = "'314.59"
$c.ClearFormats()
$ws.Range("E5").Value = "  -1.19%  "
$c = $ws.Range("D6")
$c.Value = "'98.43"
$c.ClearFormats()
$ws.Range("E6").Value = "  -7.03%  "
$ws.Range("E7").Value = "  -3.15%  "
$ws.Range("E8").Value = "  +0.21%  "
$c = $ws.Range("D9")
$c.Value = "'0.530"
$c.ClearFormats()
$ws.Range("E9").Value = "  -7.54%  "
$c = $ws.Range("D10")
$c.Value = "'35.76"
$c.ClearFormats()
$ws.Range("E10").Value = "  -8.45%  "
$c = $ws.Range("D11")
$c.Value = "'0.0819"
$c.ClearFormats()
$ws.Range("E11").Value = "  -2.79%  "
$c = $ws.Range("D12")
$c.Value = "'7.34"
$c.ClearFormats()
$ws.Range("E12").Value = "  -7.41%  "
$ws.Range("E13").Value = "  -2.63%  "
$ws.Range("D14").Value = "2.576.56"
$ws.Range("E14").Value = "  -1.98%  "
$ws.Range("D15").Value = "2.239.26"
$ws.Range("E15").Value = "  -2.31%  "
$c = $ws.Range("D16")
$c.Value = "'0.836"
$c.ClearFormats()
$ws.Range("E16").Value = "  -5.20%  "
$c = $ws.Range("D17")
$c.Value = "'13.88"
$c.ClearFormats()
$ws.Range("E17").Value = "  -5.23%  "
$ws.Range("D18").Value = "43.732.73"
$ws.Range("E18").Value = "  -0.88%  "
$c = $ws.Range("D19")
$c.Value = "'13.09"
$c.ClearFormats()
$ws.Range("E19").Value = "  -7.57%  "
$ws.Range("D20").Value = "0.0₃0968"
$ws.Range("E20").Value = "  -3.24%  "
$ws.Range("E21").Value = "  -4.31%  "
$c = $ws.Range("D22")
$c.Value = "'66.00"
$c.ClearFormats()
$ws.Range("E22").Value = "  -0.42%  "
$c = $ws.Range("D23")
$c.Value = "'235.97"
$c.ClearFormats()
$ws.Range("E23").Value = "  -0.82%  "
$c = $ws.Range("D24")
$c.Value = "'2.98"
$c.ClearFormats()
$ws.Range("E24").Value = "  -7.36%  "
$ws.Range("E25").Value = "  -8.66%  "
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("E27").Value = "  -2.22%  "
$c = $ws.Range("D28")
$c.Value = "'2.14"
$c.ClearFormats()
$ws.Range("E28").Value = "  -3.10%  "
$c = $ws.Range("D29")
$c.Value = "'36.50"
$c.ClearFormats()
$ws.Range("E29").Value = "  -6.72%  "
$c = $ws.Range("D30")
$c.Value = "'5.97"
$c.ClearFormats()
$ws.Range("E30").Value = "  -8.88%  "
$c = $ws.Range("D31")
$c.Value = "'19.98"
$c.ClearFormats()
$ws.Range("E31").Value = "  -2.78%  "
$c = $ws.Range("D32")
$c.Value = "'155.33"
$c.ClearFormats()
$ws.Range("E32").Value = "  -5.16%  "
$ws.Range("E33").Value = "  -6.41%  "
$c = $ws.Range("D34")
$c.Value = "'3.32"
$c.ClearFormats()
$ws.Range("E34").Value = "  +1.36%  "
$c = $ws.Range("D35")
$c.Value = "'2.65"
$c.ClearFormats()
$ws.Range("E35").Value = "  -3.05%  "
$c = $ws.Range("D36")
$c.Value = "'1.91"
$c.ClearFormats()
$ws.Range("E36").Value = "  -8.22%  "
$c = $ws.Range("D37")
$c.Value = "'0.108"
$c.ClearFormats()
$ws.Range("E37").Value = "  -6.75%  "
$ws.Range("E38").Value = "  -3.42%  "
$c = $ws.Range("D39")
$c.Value = "'15.59"
$c.ClearFormats()
$ws.Range("E39").Value = "  +1.01%  "
$ws.Range("E40").Value = "  -11.45%  "
$c = $ws.Range("D41")
$c.Value = "'3.99"
$c.ClearFormats()
$ws.Range("E41").Value = "  -11.43%  "
$c = $ws.Range("D42")
$c.Value = "'0.0307"
$c.ClearFormats()
$ws.Range("E42").Value = "  -6.06%  "
$ws.Range("D44").Value = "1.698.24"
$ws.Range("E44").Value = "  -3.87%  "
$c = $ws.Range("D45")
$c.Value = "'82.48"
$c.ClearFormats()
$ws.Range("E45").Value = "  -3.95%  "
$ws.Range("E46").Value = "  -6.82%  "
$ws.Range("E47").Value = "  -4.72%  "
$c = $ws.Range("D48")
$c.Value = "'101.49"
$c.ClearFormats()
$ws.Range("E48").Value = "  -2.89%  "
$c = $ws.Range("D49")
$c.Value = "'71.03"
$c.ClearFormats()
$ws.Range("E49").Value = "  -6.19%  "
$c = $ws.Range("D50")
$c.Value = "'56.11"
$c.ClearFormats()
$ws.Range("E50").Value = "  -5.98%  "
$ws.Range("E51").Value = "  -5.01%  "
